$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.33073997497559
$ws.Range("C3").Value = 17.8382396697998
$ws.Range("C4").Value = 17.62199401855469
$ws.Range("C5").Value = 17.51208305358887
$ws.Range("C6").Value = 18.20683479309082
